$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2021-04-06 "actualizacion automatica" run re-derived the fuente-y-fecha
# mapping: the 2016 entry disappeared and the 2015 entry's URL text picked up
# a run of trailing dashes (upstream artefact, reproduced verbatim here).

# Update the B1 display text / value (the hyperlink address itself is left
# pointing at the original 2015 URL).
$ws.Range("B1").Value2 = "http://opendata.aragon.es/kos/iaest/fuente-y-fecha/catastro-diciembre-de-2015--------------------------"

# Drop the hyperlink that lives on B2 before the row disappears so it doesn't
# end up orphaned in the hyperlinks collection.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.Delete()
    }
}

# Remove row 2 (Catastro diciembre de 2016 + its link) entirely.
$ws.Rows.Item(2).Delete()
